$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column at E (shifts old E:H -> F:I, and auto-updates the
#     stock data validation sqref from F to G, and the sheet dimension). ---
$ws.Columns("E").Insert()

# --- Give the new column E the same visual width as column D, and tag it
#     with the accounting-style number format (builtin id 40) used for the
#     new UnitPrice values. ---
$ws.Columns("E").ColumnWidth = $ws.Columns("D").ColumnWidth
$ws.Columns("E").NumberFormat = "#,##0.00_);[Red](#,##0.00)"

# --- Header + data for the new UnitPrice column. ---
$ws.Range("E1").Value = "UnitPrice"
$ws.Range("E2").Value = 40.05
$ws.Range("E2").NumberFormat = "#,##0.00_);[Red](#,##0.00)"

# --- Conditional-formatting ranges don't auto-shift with a column insert in
#     this engine, so fix them up by hand. Grab both FormatCondition objects
#     up front (before mutating either) since AppliesTo ranges overlap once
#     the D rule grows onto the new column. ---
$fcD = $ws.Range("D2:D1048576").FormatConditions.Item(1)
$fcE = $ws.Range("E2:E1048576").FormatConditions.Item(1)

# Old HSN_Code duplicate-check (was column E, now column F).
$fcE.ModifyAppliesToRange($ws.Range("F2:F1048576"))
# Old Item_Code duplicate-check (column D) now also covers the new column E.
$fcD.ModifyAppliesToRange($ws.Range("D2:E1048576"))

# --- Reflect where the editor's cursor ended up after the edit. ---
$ws.Range("G8").Select()
